$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values: row 2-4 updated (FAPs-only -> ECs sending), rows 5-7 added (FAPs sending, restoring original FAPs rows)
$data = @(
    @("ECs", "Wnt2", "Fzd9", "ECs", 1, 0.3333333333333333, 0.01070233333333333, 0.032107, 0.004227647500550067, 0.004227647500550067, 3, 1, 0.083686, 0.251058, 0.08174459316063268, 0.08174459316063269, 0.0008956354673333333, 0.008060719206, 0.0003455873249590308, 0.0003455873249590309),
    @("ECs", "Wnt2", "Fzd9", "FAPs", 1, 0.3333333333333333, 0.01070233333333333, 0.032107, 0.004227647500550067, 0.004227647500550067, 3, 1, 0.6105696666666667, 1.831709, 0.5964052409947874, 0.5964052409947875, 0.006534520095888889, 0.05881068086299999, 0.002521391126406573, 0.002521391126406574),
    @("ECs", "Wnt2", "Fzd9", "MuSCs", 1, 0.3333333333333333, 0.01070233333333333, 0.032107, 0.004227647500550067, 0.004227647500550067, 3, 1, 0.329494, 0.988482, 0.3218501658445798, 0.3218501658445799, 0.003526354619333333, 0.031737191574, 0.001360669049184462, 0.001360669049184463),
    @("FAPs", "Wnt2", "Fzd9", "ECs", 3, 1, 2.520808, 7.562424, 0.99577235249945, 0.99577235249945, 3, 1, 0.083686, 0.251058, 0.08174459316063268, 0.08174459316063269, 0.210956338288, 1.898607044592, 0.08139900583567365, 0.08139900583567367),
    @("FAPs", "Wnt2", "Fzd9", "FAPs", 3, 1, 2.520808, 7.562424, 0.99577235249945, 0.99577235249945, 3, 1, 0.6105696666666667, 1.831709, 0.5964052409947874, 0.5964052409947875, 1.539128900290667, 13.852160102616, 0.5938838498683808, 0.5938838498683809),
    @("FAPs", "Wnt2", "Fzd9", "MuSCs", 3, 1, 2.520808, 7.562424, 0.99577235249945, 0.99577235249945, 3, 1, 0.329494, 0.988482, 0.3218501658445798, 0.3218501658445799, 0.8305911111520001, 7.475320000368, 0.3204894967953953, 0.3204894967953954),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = 1 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
